# Applies the 25-aout-2025 contract edit:
#  1) Justify ("both") the paragraph "La Compagnie Financiere Africaine du Gabon..."
#  2) Justify ("both") the paragraph "COFINA Gabon SA est representee par..."
#  3) Swap the signatory from Monsieur El Hadji Mamadou FAYE (DG) to
#     Madame Jenny MVOU (DGA) in that same paragraph.

$d = $word.ActiveDocument

# wdAlignParagraphJustify = 3
$wdAlignParagraphJustify = 3

# --- 1): justify the "La Compagnie Financiere..." paragraph ---------------

$companyRange = $d.Content.Duplicate
$companyRange.Find.ClearFormatting()
$companyRange.Find.Text = "La Compagnie Financière Africaine du Gabon"
$foundCompany = $companyRange.Find.Execute()
if (-not $foundCompany) {
    throw "Could not find the 'La Compagnie Financiere...' paragraph"
}
$companyRange.ParagraphFormat.Alignment = $wdAlignParagraphJustify

# --- 2): justify the "COFINA Gabon SA est representee par..." paragraph ---

$signatoryRange = $d.Content.Duplicate
$signatoryRange.Find.ClearFormatting()
$signatoryRange.Find.Text = "COFINA Gabon SA est représentée par"
$foundSignatory = $signatoryRange.Find.Execute()
if (-not $foundSignatory) {
    throw "Could not find the 'COFINA Gabon SA est representee par...' paragraph"
}
$signatoryRange.ParagraphFormat.Alignment = $wdAlignParagraphJustify

# --- 3): replace the "Monsieur ... FAYE ... Directeur General" text -------

# a) "Monsieur " -> "Madame " (plain, non-bold run)
$f1 = $d.Content.Duplicate
$f1.Find.ClearFormatting()
$f1.Find.Text = "Monsieur "
$found1 = $f1.Find.Execute()
if (-not $found1) {
    throw "Could not find 'Monsieur '"
}
$f1.Text = "Madame "

# b) "El Hadji Mamadou FAYE" (bold run) -> "Jenny MVOU, "
$f2 = $d.Content.Duplicate
$f2.Find.ClearFormatting()
$f2.Find.Text = "El Hadji Mamadou FAYE"
$found2 = $f2.Find.Execute()
if (-not $found2) {
    throw "Could not find 'El Hadji Mamadou FAYE'"
}
$f2.Text = "Jenny MVOU, "

# c) ", en qualité de Directeur Général, dument habilité aux fins des présentes, "
#    -> "en qualité de Directeur Général Adjointe, dument habilitée aux fins des présentes, "
$f3 = $d.Content.Duplicate
$f3.Find.ClearFormatting()
$f3.Find.Text = ", en qualité de Directeur Général, dument habilité aux fins des présentes, "
$found3 = $f3.Find.Execute()
if (-not $found3) {
    throw "Could not find the 'en qualite de Directeur General...' sentence"
}
$f3.Text = "en qualité de Directeur Général Adjointe, dument habilitée aux fins des présentes, "
